# The deck originally ships the "Integral" theme on the slide master
# (ppt/theme/theme1.xml) and the stock "Office Theme" colours on the
# notes master (ppt/theme/theme2.xml). The authored change swaps the
# two themes' content, so the slides end up styled with the plain
# Office palette.
#
# Re-create that by pushing the Office theme's 12 scheme colours (in
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order) onto the
# presentation's live theme through the Design/Theme object model.

$p = $ppt.ActivePresentation

$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Colors($i).RGB = $officeColors[$i - 1]
}
